$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1308.5883
$ws.Range("I19").Value = 902.36365
$ws.Range("J19").Value = 2053.3333
$ws.Range("K19").Value = 902.36365
$ws.Range("L19").Value = 2053.3333
$ws.Range("M19").Value = -727.36365
$ws.Range("N19").Value = -2403.3333
$ws.Range("H82").Value = 3460.875
$ws.Range("I82").Value = 1232.3334
$ws.Range("J82").Value = 4798
$ws.Range("K82").Value = 3697.0002
$ws.Range("L82").Value = 14394
$ws.Range("M82").Value = -3291.0002
$ws.Range("N82").Value = -15206
$ws.Range("H85").Value = 3460.875
$ws.Range("I85").Value = 1232.3334
$ws.Range("J85").Value = 4798
$ws.Range("K85").Value = 3697.0002
$ws.Range("L85").Value = 14394
$ws.Range("M85").Value = -2293.0002
$ws.Range("N85").Value = -17202
$ws.Range("H132").Value = 5108058.5
$ws.Range("I132").Value = 6417056.5
$ws.Range("K132").Value = 19251169.5
$ws.Range("M132").Value = -19248639.5
$ws.Range("H138").Value = 4204.345
$ws.Range("I138").Value = 2347.7
$ws.Range("J138").Value = 5265.2856
$ws.Range("K138").Value = 7043.099999999999
$ws.Range("L138").Value = 15795.8568
$ws.Range("M138").Value = -1903.099999999999
$ws.Range("N138").Value = -26075.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 40325.684
$ws.Range("I32").Value = 18118.613
$ws.Range("J32").Value = 119438.375
$ws.Range("K32").Value = 18118.613
$ws.Range("L32").Value = 119438.375
$ws.Range("M32").Value = -17831.613
$ws.Range("N32").Value = -120012.375
$ws.Range("H45").Value = 63706.438
$ws.Range("I45").Value = 111897
$ws.Range("J45").Value = 1747.1428
$ws.Range("K45").Value = 111897
$ws.Range("L45").Value = 1747.1428
$ws.Range("M45").Value = -111520
$ws.Range("N45").Value = -2501.1428
$ws.Range("H61").Value = 2694.0715
$ws.Range("I61").Value = 2719.6667
$ws.Range("J61").Value = 2648
$ws.Range("K61").Value = 2719.6667
$ws.Range("L61").Value = 2648
$ws.Range("M61").Value = -2507.6667
$ws.Range("N61").Value = -3072
$ws.Range("H63").Value = 2800
$ws.Range("I63").Value = 1000
$ws.Range("J63").Value = 3400
$ws.Range("K63").Value = 1000
$ws.Range("L63").Value = 3400
$ws.Range("N63").Value = -4772
$ws.Range("M63").Value = -314
$ws.Range("H66").Value = 2800
$ws.Range("I66").Value = 1000
$ws.Range("J66").Value = 3400
$ws.Range("K66").Value = 5000
$ws.Range("L66").Value = 17000
$ws.Range("N66").Value = -23864
$ws.Range("M66").Value = -1568
$ws.Range("H102").Value = 73701.21000000001
$ws.Range("I102").Value = 113140.78
$ws.Range("J102").Value = 2710
$ws.Range("K102").Value = 113140.78
$ws.Range("L102").Value = 2710
$ws.Range("M102").Value = -111518.78
$ws.Range("N102").Value = -5954
$ws.Range("H117").Value = 44980
$ws.Range("J117").Value = 44980
$ws.Range("L117").Value = 44980
$ws.Range("N117").Value = -54158
$ws.Range("H136").Value = 2694.0715
$ws.Range("I136").Value = 2719.6667
$ws.Range("J136").Value = 2648
$ws.Range("K136").Value = 8159.000100000001
$ws.Range("L136").Value = 7944
$ws.Range("M136").Value = -5609.000100000001
$ws.Range("N136").Value = -13044

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 3164
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("H85").Value = 3164
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("H105").Value = 126933.125
$ws.Range("I105").Value = 93009
$ws.Range("J105").Value = 201566.2
$ws.Range("K105").Value = 93009
$ws.Range("L105").Value = 201566.2
$ws.Range("M105").Value = -91262
$ws.Range("N105").Value = -205060.2
$ws.Range("H134").Value = 3891.4934
$ws.Range("I134").Value = 3862.7258
$ws.Range("J134").Value = 4028.6924
$ws.Range("K134").Value = 11588.1774
$ws.Range("L134").Value = 12086.0772
$ws.Range("M134").Value = -9053.1774
$ws.Range("N134").Value = -17156.0772
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 855.6957
$ws.Range("I16").Value = 641.75
$ws.Range("J16").Value = 1344.7142
$ws.Range("K16").Value = 641.75
$ws.Range("L16").Value = 1344.7142
$ws.Range("M16").Value = -354.75
$ws.Range("N16").Value = -1918.7142
$ws.Range("H31").Value = 54670.83
$ws.Range("I31").Value = 2551.2
$ws.Range("J31").Value = 82102.21000000001
$ws.Range("K31").Value = 2551.2
$ws.Range("L31").Value = 82102.21000000001
$ws.Range("M31").Value = -2256.2
$ws.Range("N31").Value = -82692.21000000001
$ws.Range("H34").Value = 54670.83
$ws.Range("I34").Value = 2551.2
$ws.Range("J34").Value = 82102.21000000001
$ws.Range("K34").Value = 2551.2
$ws.Range("L34").Value = 82102.21000000001
$ws.Range("M34").Value = -2349.2
$ws.Range("N34").Value = -82506.21000000001
$ws.Range("H43").Value = 27657
$ws.Range("J43").Value = 27657
$ws.Range("L43").Value = 27657
$ws.Range("N43").Value = -28025
$ws.Range("H58").Value = 2203.88
$ws.Range("I58").Value = 2405.5833
$ws.Range("J58").Value = 2017.6923
$ws.Range("K58").Value = 2405.5833
$ws.Range("L58").Value = 2017.6923
$ws.Range("M58").Value = -2202.5833
$ws.Range("N58").Value = -2423.6923
$ws.Range("H101").Value = 27657
$ws.Range("J101").Value = 27657
$ws.Range("L101").Value = 27657
$ws.Range("N101").Value = -34147
$ws.Range("H113").Value = 855.6957
$ws.Range("I113").Value = 641.75
$ws.Range("J113").Value = 1344.7142
$ws.Range("K113").Value = 641.75
$ws.Range("L113").Value = 1344.7142
$ws.Range("M113").Value = 1528.25
$ws.Range("N113").Value = -5684.7142
$ws.Range("H131").Value = 36833.332
$ws.Range("J131").Value = 36833.332
$ws.Range("L131").Value = 36833.332
$ws.Range("N131").Value = -46913.332
$ws.Range("H132").Value = 3888.3125
$ws.Range("I132").Value = 4428.5454
$ws.Range("K132").Value = 13285.6362
$ws.Range("M132").Value = -10755.6362
$ws.Range("H136").Value = 2203.88
$ws.Range("I136").Value = 2405.5833
$ws.Range("J136").Value = 2017.6923
$ws.Range("K136").Value = 7216.749899999999
$ws.Range("L136").Value = 6053.0769
$ws.Range("M136").Value = -4666.749899999999
$ws.Range("N136").Value = -11153.0769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1950.4762
$ws.Range("J9").Value = 1950.4762
$ws.Range("L9").Value = 5851.4286
$ws.Range("N9").Value = -6299.4286
$ws.Range("H19").Value = 995
$ws.Range("J19").Value = 995
$ws.Range("L19").Value = 2985
$ws.Range("N19").Value = -3333
$ws.Range("H20").Value = 933.3333
$ws.Range("J20").Value = 600
$ws.Range("L20").Value = 1800
$ws.Range("N20").Value = -2254
$ws.Range("H22").Value = 3044.2092
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 3143.9268
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 9431.7804
$ws.Range("N22").Value = -9769.7804
$ws.Range("M22").Value = -2831
$ws.Range("H26").Value = 1211.0714
$ws.Range("I26").Value = 221.57143
$ws.Range("J26").Value = 2200.5715
$ws.Range("K26").Value = 664.71429
$ws.Range("L26").Value = 6601.7145
$ws.Range("M26").Value = -376.71429
$ws.Range("N26").Value = -7177.7145
$ws.Range("H27").Value = 3044.2092
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 3143.9268
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 9431.7804
$ws.Range("N27").Value = -9635.7804
$ws.Range("M27").Value = -2898
$ws.Range("H38").Value = 164.5
$ws.Range("J38").Value = 194.25
$ws.Range("L38").Value = 582.75
$ws.Range("N38").Value = -1276.75
$ws.Range("H131").Value = 724393.9
$ws.Range("I131").Value = 750
$ws.Range("J131").Value = 768251.0600000001
$ws.Range("K131").Value = 2250
$ws.Range("L131").Value = 2304753.18
$ws.Range("M131").Value = 2790
$ws.Range("N131").Value = -2314833.18

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("H122").Value = 2584.5625
$ws.Range("J122").Value = 2724
$ws.Range("L122").Value = 8172
$ws.Range("N122").Value = -13072
$ws.Range("H132").Value = 3734.7334
$ws.Range("I132").Value = 2796.3684
$ws.Range("K132").Value = 8389.1052
$ws.Range("M132").Value = -5859.1052
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5076.8184
$ws.Range("I46").Value = 9340
$ws.Range("J46").Value = 4129.4443
$ws.Range("K46").Value = 9340
$ws.Range("L46").Value = 4129.4443
$ws.Range("M46").Value = -9152
$ws.Range("N46").Value = -4505.4443
$ws.Range("H122").Value = 5722.4
$ws.Range("I122").Value = 3998.8
$ws.Range("J122").Value = 7446
$ws.Range("K122").Value = 11996.4
$ws.Range("L122").Value = 22338
$ws.Range("M122").Value = -9546.400000000001
$ws.Range("N122").Value = -27238
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9923.076999999999
$ws.Range("I15").Value = 9000
$ws.Range("K15").Value = 9000
$ws.Range("M15").Value = -8712
$ws.Range("H54").Value = 6806
$ws.Range("J54").Value = 7051.3335
$ws.Range("L54").Value = 7051.3335
$ws.Range("N54").Value = -8091.3335
$ws.Range("H132").Value = 18917.684
$ws.Range("I132").Value = 2403.8096
$ws.Range("J132").Value = 57450.055
$ws.Range("K132").Value = 7211.4288
$ws.Range("L132").Value = 172350.165
$ws.Range("M132").Value = -4681.4288
$ws.Range("N132").Value = -177410.165
